# Updates crypto price/volume figures to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    # Force the cell to stay text-typed (matches the original inlineStr
    # cells) instead of Excel auto-converting numeric-looking strings to
    # real numbers.
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Column D (Price) updates
Set-TextValue "D2"  "242.69"
Set-TextValue "D3"  "23.12"
Set-TextValue "D4"  "5.383"
Set-TextValue "D5"  "0.05980"
Set-TextValue "D6"  "3.402"
Set-TextValue "D7"  "6.485"
Set-TextValue "D8"  "0.8128"
Set-TextValue "D9"  "0.9097"
Set-TextValue "D10" "0.1411"
Set-TextValue "D11" "0.07412"
Set-TextValue "D12" "0.03354"
Set-TextValue "D13" "0.03059"
Set-TextValue "D14" "0.09335"
Set-TextValue "D15" "3.856"
Set-TextValue "D16" "0.001573"
Set-TextValue "D17" "0.04636"
Set-TextValue "D19" "0.006091"
Set-TextValue "D20" "0.005020"
Set-TextValue "D21" "0.0009879"
Set-TextValue "D22" "0.00007799"
Set-TextValue "D24" "3.615"
Set-TextValue "D40" "0.03888"
Set-TextValue "D41" "0.006206"
Set-TextValue "D44" "0.007207"
Set-TextValue "D45" "0.00005191"
Set-TextValue "D49" "0.002297"

# Column E (Volume(1h)) updates
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"
